$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns AD/AE/AF ("Wins", "Losses", "Ties") added right after the
# existing "Unnamed: 28" column (AC) - team record appended to the roster.

# Pick up the header formatting (bold font, borders, center/top alignment)
# from the existing last header cell (AC1) and stamp it onto the three new
# header cells before writing their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-45) gets the team's season record: 95-67-0.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 95
    $ws.Cells.Item($r, 31).Value = 67
    $ws.Cells.Item($r, 32).Value = 0
}
